$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a plain text value to a cell without Excel coercing it into a
# number/date. We prefix with an apostrophe (forces text entry, same as
# typing it manually) and then reset the resulting style back to "Normal" so
# we don't leave a stray quote-prefixed / custom-number-format style behind.
function Set-TextValue($addr, $text) {
    $ws.Range($addr).Value = "'" + $text
    $ws.Range($addr).Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "68.403.74"
$ws.Range("E2").Value = "  +1.18%  "

# Row 3 - Ethereum
Set-TextValue "D3" "3.751.41"
$ws.Range("E3").Value = "  -0.68%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.13%  "

# Row 5 - BNB
Set-TextValue "D5" "595.65"

# Row 6 - Solana
Set-TextValue "D6" "167.00"
$ws.Range("E6").Value = "  -1.16%  "

# Row 7 - LidoStakedEther
Set-TextValue "D7" "3.748.52"
$ws.Range("E7").Value = "  -0.74%  "

# Row 9 - XRP
$ws.Range("E9").Value = "  -0.69%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  -2.93%  "

# Row 11 - Toncoin
$ws.Range("E11").Value = "  +0.52%  "

# Row 12 - Cardano
Set-TextValue "D12" "0.448"
$ws.Range("E12").Value = "  -1.07%  "

# Row 13 - ShibaInu
$ws.Range("E13").Value = "  -5.53%  "

# Row 14 - Avalanche
Set-TextValue "D14" "36.10"
$ws.Range("E14").Value = "  -0.99%  "

# Row 15 - WrappedliquidstakedEther2.0
Set-TextValue "D15" "4.381.29"
$ws.Range("E15").Value = "  -0.61%  "

# Row 16 - WrappedEther
Set-TextValue "D16" "3.746.71"
$ws.Range("E16").Value = "  -0.70%  "

# Row 17 - WrappedBTC
Set-TextValue "D17" "68.447.40"
$ws.Range("E17").Value = "  +1.27%  "

# Row 18 - Chainlink
Set-TextValue "D18" "17.88"
$ws.Range("E18").Value = "  -3.92%  "

# Row 19 - Polkadot
$ws.Range("E19").Value = "  -2.19%  "

# Row 20 - TRON
$ws.Range("E20").Value = "  -0.15%  "

# Row 21 - Uniswap
Set-TextValue "D21" "10.72"
$ws.Range("E21").Value = "  +1.37%  "

# Row 22 - BitcoinCash
Set-TextValue "D22" "467.83"
$ws.Range("E22").Value = "  +0.09%  "

# Row 23 - Polygon
$ws.Range("E23").Value = "  -2.65%  "

# Row 24 - Litecoin
Set-TextValue "D24" "84.21"
$ws.Range("E24").Value = "  +0.79%  "

# Row 25 - PEPE
$ws.Range("E25").Value = "  -1.78%  "

# Row 26 - Fetch.AI
Set-TextValue "D26" "2.19"
$ws.Range("E26").Value = "  -0.37%  "

# Row 27 - InternetComputer(DFINITY)
Set-TextValue "D27" "12.04"
$ws.Range("E27").Value = "  -0.83%  "

# Row 28 - RenderToken
Set-TextValue "D28" "10.13"
$ws.Range("E28").Value = "  -1.39%  "

# Row 29 - Dai
$ws.Range("E29").Value = "  -0.16%  "

# Row 30 - WrappedeETH
Set-TextValue "D30" "3.898.13"
$ws.Range("E30").Value = "  -0.48%  "

# Row 31 - PancakeSwap
$ws.Range("E31").Value = "  -4.47%  "

# Row 32 - NEARProtocol
Set-TextValue "D32" "7.28"
$ws.Range("E32").Value = "  -4.43%  "

# Row 33 - EthereumClassic
Set-TextValue "D33" "29.87"
$ws.Range("E33").Value = "  -2.03%  "

# Row 34 - ImmutableX
$ws.Range("E34").Value = "  -1.86%  "

# Row 35 - Aptos
$ws.Range("E35").Value = "  +1.47%  "

# Row 37 - RenzoRestakedETH
Set-TextValue "D37" "3.707.14"
$ws.Range("E37").Value = "  -0.78%  "

# Row 38 - Hedera
$ws.Range("E38").Value = "  -2.24%  "

# Row 39 - dogwifhat
Set-TextValue "D39" "3.38"
$ws.Range("E39").Value = "  -10.93%  "

# Row 40 - Kaspa
Set-TextValue "D40" "0.139"
$ws.Range("E40").Value = "  +0.56%  "

# Row 41 - Mantle
Set-TextValue "D41" "0.998"
$ws.Range("E41").Value = "  -0.59%  "

# Row 42 - Filecoin
Set-TextValue "D42" "5.81"
$ws.Range("E42").Value = "  +0.32%  "

# Row 43 - FirstDigitalUSD
Set-TextValue "D43" "1.00"
$ws.Range("E43").Value = "  +0.13%  "

# Row 45 - TheGraph
$ws.Range("E45").Value = "  -1.69%  "

# Row 46 - Cosmos
Set-TextValue "D46" "8.60"
$ws.Range("E46").Value = "  -1.15%  "

# Row 47 - now Stacks (was Arweave)
$ws.Range("B47").Value = "Stacks"
$ws.Range("C47").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue "D47" "1.93"
$ws.Range("E47").Value = "  -0.29%  "

# Row 48 - now Arweave (was Stacks)
$ws.Range("B48").Value = "Arweave"
$ws.Range("C48").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
Set-TextValue "D48" "43.07"
$ws.Range("E48").Value = "  +10.72%  "

# Row 49 - OKB
Set-TextValue "D49" "45.78"
$ws.Range("E49").Value = "  -0.10%  "

# Row 50 - Monero
Set-TextValue "D50" "146.43"
$ws.Range("E50").Value = "  +4.90%  "

# Row 51 - Bittensor
Set-TextValue "D51" "391.53"
$ws.Range("E51").Value = "  -1.16%  "
